# Update status text "Ready for handoff" -> "In Translation" across all sheets,
# then resize the affected columns to match the recalculated (auto-fit) widths.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: status text lives in columns E (zh-cn) and F (de-de), rows 2-4
$overview.Range("E2:F4").Value = "In Translation"

# zh-cn / de-de sheets: status text lives in column C ("Status"), rows 2-4
$zhcn.Range("C2:C4").Value = "In Translation"
$dede.Range("C2:C4").Value = "In Translation"

# Recompute the width of the columns that held the text that changed
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
